$d = $word.ActiveDocument

# Locate the "Add Last ..." paragraph and the ".Net BLC ..." paragraph that
# follows it a few paragraphs later; together they bound the block that
# needs to be rewritten (bookmark relocation, proofErr wrapping, new run).
$pStart = $null
$pEnd = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "Add Last*") {
        $pStart = $p
    }
    if ($t -like "*LinkList*") {
        $pEnd = $p
    }
}

$start = $pStart.Range.Start
$end = $pEnd.Range.End
$rng = $d.Range($start, $end)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Add Last – you shift the pointer of the tail</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Remove last – set last pointer to null and move the tail</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">.Net BLC has a built in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LinkList</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0DF"/></w:r><w:r><w:t xml:space="preserve"> use this</w:t></w:r><w:r><w:t xml:space="preserve"> – this is a Doubly-linked Circular List</w:t></w:r></w:p>'

$rng.InsertXML($xml)
